# Daily attendance processing - 2026-01-26 14:43:57
# Swap the order of recorder names in column G ("Recorded By") so that
# "System, dnasr281@gmail.com" becomes "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
